$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.176837485174929
$ws.Range("D2").Value = 0.002974535826233637
$ws.Range("E2").Value = 0.8824200776671347
$ws.Range("F2").Value = 0.7722284337438623
$ws.Range("G2").Value = 0.6574680208752284
$ws.Range("H2").Value = 0.6343096553614487
$ws.Range("I2").Value = 0.8676753838507025
$ws.Range("L2").Value = 0.3836018548599327
$ws.Range("M2").Value = 0.3401644527974668

$ws.Range("B3").Value = 1.081232832122282
$ws.Range("D3").Value = 0.003227087922964222
$ws.Range("E3").Value = 0.8130003696447687
$ws.Range("F3").Value = 0.736615759759303
$ws.Range("G3").Value = 0.6173742649015708
$ws.Range("H3").Value = 0.6229143106936306
$ws.Range("I3").Value = 0.8786687051774038
$ws.Range("L3").Value = 0.3455242620919705
$ws.Range("M3").Value = 0.3100154504194776

$ws.Range("B4").Value = 1.022494743027039
$ws.Range("D4").Value = 0.003394258971711195
$ws.Range("E4").Value = 0.7702308310097692
$ws.Range("F4").Value = 0.7155945718421748
$ws.Range("G4").Value = 0.5935436061577377
$ws.Range("H4").Value = 0.6165600582999389
$ws.Range("I4").Value = 0.8862462179699264
$ws.Range("L4").Value = 0.3221627195844121
$ws.Range("M4").Value = 0.2915050746212913

$ws.Range("B5").Value = 0.9985505296633335
$ws.Range("D5").Value = 0.003465379380220623
$ws.Range("E5").Value = 0.7527664579437214
$ws.Range("F5").Value = 0.7072382074766352
$ws.Range("G5").Value = 0.5840272372547588
$ws.Range("H5").Value = 0.6141309166064985
$ws.Range("I5").Value = 0.8895421112252677
$ws.Range("L5").Value = 0.3126474890053146
$ws.Range("M5").Value = 0.28396257117857

$ws.Range("B6").Value = 0.9945741661600778
$ws.Range("D6").Value = 0.003477368517840307
$ws.Range("E6").Value = 0.7498644063622493
$ws.Range("F6").Value = 0.7058632443069115
$ws.Range("G6").Value = 0.5824587276371318
$ws.Range("H6").Value = 0.6137372012783544
$ws.Range("I6").Value = 0.8901019514634783
$ws.Range("L6").Value = 0.3110677874633723
$ws.Range("M6").Value = 0.282710191438504

$ws.Range("B7").Value = 1.022171853636564
$ws.Range("D7").Value = 0.003395206052935307
$ws.Range("E7").Value = 0.7699954422846105
$ws.Range("F7").Value = 0.7154810283570043
$ws.Range("G7").Value = 0.5934144803574668
$ws.Range("H7").Value = 0.6165266508606351
$ws.Range("I7").Value = 0.8862898255076175
$ws.Range("L7").Value = 0.3220343742026444
$ws.Range("M7").Value = 0.2914033508293272

$ws.Range("B8").Value = 1.143881245403861
$ws.Range("D8").Value = 0.003059075711685066
$ws.Range("E8").Value = 0.8585150423641181
$ws.Range("F8").Value = 0.759772209230448
$ws.Range("G8").Value = 0.6434784575747869
$ws.Range("H8").Value = 0.6302463907991012
$ws.Range("I8").Value = 0.8712940509485776
$ws.Range("L8").Value = 0.3704690157366315
$ws.Range("M8").Value = 0.3297689974708788

$ws.Range("B9").Value = 1.382224387801045
$ws.Range("D9").Value = 0.002498048881701687
$ws.Range("E9").Value = 1.030904948253522
$ws.Range("F9").Value = 0.8534545117699963
$ws.Range("G9").Value = 0.7480440086154374
$ws.Range("H9").Value = 0.6623102711210151
$ws.Range("I9").Value = 0.8484598970722601
$ws.Range("L9").Value = 0.4655907249104416
$ws.Range("M9").Value = 0.405004059385746

$ws.Range("B10").Value = 1.557101345408796
$ws.Range("D10").Value = 0.002148565239374012
$ws.Range("E10").Value = 1.156787668031228
$ws.Range("F10").Value = 0.9266213429432355
$ws.Range("G10").Value = 0.8289729631373746
$ws.Range("H10").Value = 0.6891014619245368
$ws.Range("I10").Value = 0.8357009936032753
$ws.Range("L10").Value = 0.5355664369157864
$ws.Range("M10").Value = 0.4602725725576278

$ws.Range("B11").Value = 1.636601298962319
$ws.Range("D11").Value = 0.002003821386211246
$ws.Range("E11").Value = 1.213879516732646
$ws.Range("F11").Value = 0.9608875817463769
$ws.Range("G11").Value = 0.8667266878083524
$ws.Range("H11").Value = 0.7020111113183134
$ws.Range("I11").Value = 0.8307716699318846
$ws.Range("L11").Value = 0.567421140717471
$ws.Range("M11").Value = 0.4854133007046215

$ws.Range("B12").Value = 1.666697519306524
$ws.Range("D12").Value = 0.00195111123403291
$ws.Range("E12").Value = 1.235472967010793
$ws.Range("F12").Value = 0.9740075788648568
$ws.Range("G12").Value = 0.8811616136073894
$ws.Range("H12").Value = 0.7070050404674078
$ws.Range("I12").Value = 0.8290310910036496
$ws.Range("L12").Value = 0.5794868812524498
$ws.Range("M12").Value = 0.4949330787767252

$ws.Range("B13").Value = 1.660216163961195
$ws.Range("D13").Value = 0.001962368991813346
$ws.Range("E13").Value = 1.230823605101904
$ws.Range("F13").Value = 0.9711754972379367
$ws.Range("G13").Value = 0.8780465835014013
$ws.Range("H13").Value = 0.7059248004967174
$ws.Range("I13").Value = 0.8294003454573229
$ws.Range("L13").Value = 0.576888174752213
$ws.Range("M13").Value = 0.492882849783598

$ws.Range("B14").Value = 1.639077513650648
$ws.Range("D14").Value = 0.001999442468646917
$ws.Range("E14").Value = 1.215656550275128
$ws.Range("F14").Value = 0.9619640670010625
$ws.Range("G14").Value = 0.8679114642167178
$ws.Range("H14").Value = 0.7024198455152941
$ws.Range("I14").Value = 0.8306259437056553
$ws.Range("L14").Value = 0.5684137361194246
$ws.Range("M14").Value = 0.4861965096686873

$ws.Range("B15").Value = 1.626128315327264
$ws.Range("D15").Value = 0.002022426290884782
$ws.Range("E15").Value = 1.206362872699373
$ws.Range("F15").Value = 0.9563406536524894
$ws.Range("G15").Value = 0.8617215380087657
$ws.Range("H15").Value = 0.7002867197851401
$ws.Range("I15").Value = 0.8313930812814831
$ws.Range("L15").Value = 0.5632232909821937
$ws.Range("M15").Value = 0.4821008633925317

$ws.Range("B16").Value = 1.551904682426539
$ws.Range("D16").Value = 0.002158315677245337
$ws.Range("E16").Value = 1.153053005509804
$ws.Range("F16").Value = 0.9244019740254856
$ws.Range("G16").Value = 0.8265248499999132
$ws.Range("H16").Value = 0.6882724441562118
$ws.Range("I16").Value = 0.8360407607527804
$ws.Range("L16").Value = 0.5334850934228541
$ws.Range("M16").Value = 0.4586295170751669

$ws.Range("B17").Value = 1.506356628668868
$ws.Range("D17").Value = 0.002245364156236551
$ws.Range("E17").Value = 1.120304042213633
$ws.Range("F17").Value = 0.905062245166846
$ws.Range("G17").Value = 0.8051756801951342
$ws.Range("H17").Value = 0.6810880924199978
$ws.Range("I17").Value = 0.8391161887831657
$ws.Range("L17").Value = 0.5152472880365337
$ws.Range("M17").Value = 0.4442300848152314

$ws.Range("B18").Value = 1.480153722074874
$ws.Range("D18").Value = 0.002296769027932077
$ws.Range("E18").Value = 1.101451503442831
$ws.Range("F18").Value = 0.8940308287354526
$ws.Range("G18").Value = 0.7929843546874338
$ws.Range("H18").Value = 0.6770237315176928
$ws.Range("I18").Value = 0.8409674142103754
$ws.Range("L18").Value = 0.504759502385383
$ws.Range("M18").Value = 0.4359478168682855

$ws.Range("B19").Value = 1.471281066995573
$ws.Range("D19").Value = 0.002314401818316325
$ws.Range("E19").Value = 1.095065610369318
$ws.Range("F19").Value = 0.8903115436438469
$ws.Range("G19").Value = 0.7888716178085815
$ws.Range("H19").Value = 0.6756592239411248
$ws.Range("I19").Value = 0.8416083383478252
$ws.Range("L19").Value = 0.5012088873347125
$ws.Range("M19").Value = 0.4331435761200595

$ws.Range("B20").Value = 1.511205809249986
$ws.Range("D20").Value = 0.002235958946774641
$ws.Range("E20").Value = 1.123791910377918
$ws.Range("F20").Value = 0.9071114205905531
$ws.Range("G20").Value = 0.8074391876504023
$ws.Range("H20").Value = 0.6818458421304285
$ws.Range("I20").Value = 0.8387802825346355
$ws.Range("L20").Value = 0.5171885156273675
$ws.Range("M20").Value = 0.4457629410213997

$ws.Range("B21").Value = 1.645286692358354
$ws.Range("D21").Value = 0.001988495631476184
$ws.Range("E21").Value = 1.22011219635624
$ws.Range("F21").Value = 0.9646657531184246
$ws.Range("G21").Value = 0.8708846106375745
$ws.Range("H21").Value = 0.703446465449872
$ws.Range("I21").Value = 0.8302625327129647
$ws.Range("L21").Value = 0.5709028025008251
$ws.Range("M21").Value = 0.4881604649229274

$ws.Range("B22").Value = 1.732865160856818
$ws.Range("D22").Value = 0.001839036951875705
$ws.Range("E22").Value = 1.282911049276379
$ws.Range("F22").Value = 1.003122203538453
$ws.Range("G22").Value = 0.9131582265077043
$ws.Range("H22").Value = 0.7181783305693727
$ws.Range("I22").Value = 0.8254305325561972
$ws.Range("L22").Value = 0.6060259362247393
$ws.Range("M22").Value = 0.5158668266245741

$ws.Range("B23").Value = 1.686127932505485
$ws.Range("D23").Value = 0.001917664925835449
$ws.Range("E23").Value = 1.249408399648615
$ws.Range("F23").Value = 0.9825193337446336
$ws.Range("G23").Value = 0.890520899439025
$ws.Range("H23").Value = 0.7102589296521558
$ws.Range("I23").Value = 0.8279421340679036
$ws.Range("L23").Value = 0.5872784868096801
$ws.Range("M23").Value = 0.501079778565682

$ws.Range("B24").Value = 1.509013546712424
$ws.Range("D24").Value = 0.002240206811675138
$ws.Range("E24").Value = 1.122215122102034
$ws.Range("F24").Value = 0.9061847168295571
$ws.Range("G24").Value = 0.806415598918278
$ws.Range("H24").Value = 0.6815030579646191
$ws.Range("I24").Value = 0.8389318868659004
$ws.Range("L24").Value = 0.5163108948070487
$ws.Range("M24").Value = 0.4450699486763057

$ws.Range("B25").Value = 1.317785018602137
$ws.Range("D25").Value = 0.002639021062870173
$ws.Range("E25").Value = 0.9844015314295547
$ws.Range("F25").Value = 0.8273619042632134
$ws.Range("G25").Value = 0.719050622898294
$ws.Range("H25").Value = 0.6530746801680891
$ws.Range("I25").Value = 0.8539324794135013
$ws.Range("L25").Value = 0.4398423147712549
$ws.Range("M25").Value = 0.384651704609837
